# Add a new analysis column (L): "Lower quartile for population density"
# This flags countries whose Population Density (column G) is below 52,
# i.e. below the lower quartile threshold used in the analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in L1, matching the existing header row styling/behaviour.
$ws.Range("L1").Value = "Lower quartile for population density"

# Data rows 2-32 hold one country each; add the boolean-flag formula
# for every row, referencing that row's Population Density cell (G).
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 12).Formula = "=IF(G$r<52, ""True"",""False"")"
}

# Leave the selection where the author left it when finishing the edit.
$ws.Range("N24").Select() | Out-Null
